# modulo upload agentes por excel se permite subir sucursales.
# Adds three new header columns (DIRECCION, CELULAR, SUCURSALES) to the
# "Agentes" upload-template sheet, re-styling the header row (white text on
# a solid blue fill) and widening the columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells ------------------------------------------------
$ws.Range("D1").Value = "DIRECCION"
$ws.Range("E1").Value = "CELULAR"
$ws.Range("F1").Value = "SUCURSALES (cod sucursal separado por ,)"

# Give the new headers the same look (font/fill) as the existing "CORREO"
# header (C1) by copying its formatting across.
$ws.Range("C1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Header row formatting: white text on solid blue fill -----------
# NOTE: Excel COM colors are packed as 0x00BBGGRR, so RGB(0x00,0x70,0xC0)
# ("#0070C0") is written as 0xC07000.
$header = $ws.Range("A1:F1")
$header.Font.ThemeColor = 2
$header.Interior.Color = 0xC07000

# --- 3. Column widths so the new data is readable -----------------------
$ws.Columns.Item(2).ColumnWidth = 18.877604166666668
$ws.Columns.Item(3).ColumnWidth = 15.022135416666666
$ws.Columns.Item(4).ColumnWidth = 10.736979166666666
$ws.Columns.Item(5).ColumnWidth = 10.736979166666666
$ws.Columns.Item(6).ColumnWidth = 37.451822916666664

Write-Host "Agentes template updated: added DIRECCION, CELULAR, SUCURSALES columns"
